$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (37) down to the new row (38)
# so the new row matches the existing style (date format, etc.) exactly.
$ws.Range("A37:C37").Copy()
$ws.Range("A38:C38").PasteSpecial(-4122)

# Fill in the new time-sheet entry: 4/6/2010, 1 hour, Weekly Meeting
$ws.Range("A38").Value2 = 40274
$ws.Range("B38").Value2 = 1
$ws.Range("C38").Value = "Weekly Meeting"

# Move the active selection down, same as Excel would after entering a row
$ws.Range("A39").Select()
